$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.67
$ws.Range("H2").Value = 6.2
$ws.Range("I2").Value = 7.2
$ws.Range("Q2").Value = 1.96
$ws.Range("S2").Value = 3.5
$ws.Range("T2").Value = 1.97
$ws.Range("U2").Value = 1.85
$ws.Range("V2").Value = 1.16
$ws.Range("W2").Value = 2.48
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 22
$ws.Range("Z2").Value = 60
$ws.Range("AA2").Value = 220
$ws.Range("AB2").Value = 8.2
$ws.Range("AD2").Value = 27
$ws.Range("AE2").Value = 130
$ws.Range("AH2").Value = 25
$ws.Range("AI2").Value = 120
$ws.Range("AJ2").Value = 16.5
$ws.Range("AM2").Value = 170
$ws.Range("AO2").Value = 170

# Row 3
$ws.Range("P3").Value = 3.6
$ws.Range("Q3").Value = 1.29
$ws.Range("R3").Value = 2.06
$ws.Range("S3").Value = 1.72
$ws.Range("T3").Value = 1.74
$ws.Range("Y3").Value = 65
$ws.Range("AD3").Value = 970
$ws.Range("AF3").Value = 12.5
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 14

# Row 4
$ws.Range("I4").Value = 2.28
$ws.Range("J4").Value = 3.3
$ws.Range("K4").Value = 3.35
$ws.Range("T4").Value = 2.06
$ws.Range("V4").Value = 1.78
$ws.Range("W4").Value = 1.33
$ws.Range("X4").Value = 9.6
$ws.Range("Z4").Value = 12
$ws.Range("AG4").Value = 16.5
$ws.Range("AJ4").Value = 80
$ws.Range("AK4").Value = 60
$ws.Range("AM4").Value = 150

# Row 5
$ws.Range("O5").Value = 1.44
$ws.Range("Q5").Value = 2.3
$ws.Range("T5").Value = 1.95
$ws.Range("AI5").Value = 55
$ws.Range("AM5").Value = 120

# Row 7
$ws.Range("F7").Value = 1.49
$ws.Range("G7").Value = 1.84
$ws.Range("I7").Value = 9.2
$ws.Range("K7").Value = 7.8
$ws.Range("L7").Value = 1.32
$ws.Range("N7").Value = 1.87
$ws.Range("P7").Value = 1.87
$ws.Range("Q7").Value = 1.91
$ws.Range("R7").Value = 1.16
$ws.Range("S7").Value = 1.91
$ws.Range("W7").Value = 2.18
